$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# Header text updates (rich-text runs share identical formatting, so a
# plain Value assignment reproduces the same visible text)
$ws.Range("A8").Value = "Volume 32   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# Weekly crime statistics data updates.
# A few cells switch between the "no data" text placeholders ("0" /
# "***.*") and real numbers, which also changes their number format.
# Copying format from a same-shaped, never-edited reference cell (row 14)
# keeps the style table untouched while Value assigns the real content.
$ws.Range("I14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -26.666666666666
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 57.142857142857
$ws.Range("N15").Value = -42.105263157894
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -7.692307692307
$ws.Range("I16").Value = 74
$ws.Range("J16").Value = 97
$ws.Range("K16").Value = -23.711340206185
$ws.Range("L16").Value = -16.853932584269
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -87.888707037643
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 169
$ws.Range("J17").Value = 198
$ws.Range("K17").Value = -14.646464646464
$ws.Range("L17").Value = -19.138755980861
$ws.Range("M17").Value = 38.524590163934
$ws.Range("N17").Value = -11.979166666666
$ws.Range("C18").Value = 4
$ws.Range("D14").Copy($ws.Range("D18"))
$ws.Range("D18").Value = "0"
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = "***.*"
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 12
$ws.Range("L18").Value = -34.117647058823
$ws.Range("M18").Value = -66.666666666666
$ws.Range("N18").Value = -92.612137203166
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -54.545454545454
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 169
$ws.Range("J19").Value = 188
$ws.Range("K19").Value = -10.106382978723
$ws.Range("L19").Value = -5.05617977528
$ws.Range("M19").Value = -8.648648648648
$ws.Range("N19").Value = -42.905405405405
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 15
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 38
$ws.Range("G20").Value = 38
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 177
$ws.Range("J20").Value = 163
$ws.Range("K20").Value = 8.588957055214
$ws.Range("L20").Value = 21.232876712328
$ws.Range("M20").Value = 25.531914893617
$ws.Range("N20").Value = -90.737833594976
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -17.142857142857
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 120
$ws.Range("H21").Value = -6.666666666666
$ws.Range("I21").Value = 658
$ws.Range("J21").Value = 714
$ws.Range("K21").Value = -7.843137254901
$ws.Range("L21").Value = -8.864265927977
$ws.Range("M21").Value = -14.877102199223
$ws.Range("N21").Value = -82.643102083882
$ws.Range("I14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("D14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = "0"
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = "***.*"
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = 20
$ws.Range("L22").Value = -14.285714285714
$ws.Range("M22").Value = -33.333333333333
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -15.384615384615
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 92
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 692
$ws.Range("J24").Value = 779
$ws.Range("K24").Value = -11.168164313222
$ws.Range("L24").Value = -16.525934861278
$ws.Range("M24").Value = 57.630979498861
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -27.272727272727
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -40.816326530612
$ws.Range("I25").Value = 336
$ws.Range("J25").Value = 366
$ws.Range("K25").Value = -8.196721311475
$ws.Range("L25").Value = 24.444444444444
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 64
$ws.Range("H26").Value = -20.3125
$ws.Range("I26").Value = 316
$ws.Range("J26").Value = 353
$ws.Range("K26").Value = -10.481586402266
$ws.Range("L26").Value = 7.482993197278
$ws.Range("M26").Value = -8.93371757925
$ws.Range("I14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("I27").Value = 18
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = -14.285714285714
$ws.Range("C28").Value = 1
$ws.Range("I14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 2
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 233.333333333333
$ws.Range("I28").Value = 47
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = 104.347826086957
$ws.Range("L28").Value = 27.027027027027
$ws.Range("I33").Value = 4
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 300
